$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 894654.0600000001
$ws.Range("I15").Value = 894654.0600000001
$ws.Range("K15").Value = 2683962.18
$ws.Range("M15").Value = -2683793.18
$ws.Range("H17").Value = 1634.65
$ws.Range("J17").Value = 1634.65
$ws.Range("L17").Value = 4903.950000000001
$ws.Range("N17").Value = -5239.950000000001
$ws.Range("H119").Value = 500
$ws.Range("J119").Value = 500
$ws.Range("L119").Value = 1500
$ws.Range("N119").Value = -11176
$ws.Range("H125").Value = 399999.38
$ws.Range("J125").Value = 453570.72
$ws.Range("L125").Value = 4082136.48
$ws.Range("N125").Value = -4087056.48
$ws.Range("H132").Value = 3644.4688
$ws.Range("I132").Value = 3504.2415
$ws.Range("K132").Value = 10512.7245
$ws.Range("M132").Value = -7982.7245
$ws.Range("H133").Value = 69894.60000000001
$ws.Range("J133").Value = 69894.60000000001
$ws.Range("L133").Value = 69894.60000000001
$ws.Range("N133").Value = -80014.60000000001
$ws.Range("H137").Value = 3298.8215
$ws.Range("J137").Value = 3925.842
$ws.Range("L137").Value = 11777.526
$ws.Range("N137").Value = -16877.526
$ws.Range("H138").Value = 5255.154
$ws.Range("J138").Value = 7060
$ws.Range("L138").Value = 21180
$ws.Range("N138").Value = -31460

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1949.1666
$ws.Range("I2").Value = 2360.6365
$ws.Range("J2").Value = 1302.5714
$ws.Range("K2").Value = 2360.6365
$ws.Range("L2").Value = 1302.5714
$ws.Range("M2").Value = -2247.6365
$ws.Range("N2").Value = -1528.5714
$ws.Range("H61").Value = 4687.048
$ws.Range("I61").Value = 1732.1177
$ws.Range("K61").Value = 1732.1177
$ws.Range("M61").Value = -1520.1177
$ws.Range("H74").Value = 6406.952
$ws.Range("J74").Value = 10977.889
$ws.Range("L74").Value = 10977.889
$ws.Range("N74").Value = -12725.889
$ws.Range("H77").Value = 6406.952
$ws.Range("J77").Value = 10977.889
$ws.Range("L77").Value = 54889.44499999999
$ws.Range("N77").Value = -63625.44499999999
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("H116").Value = 1949.1666
$ws.Range("I116").Value = 2360.6365
$ws.Range("J116").Value = 1302.5714
$ws.Range("K116").Value = 2360.6365
$ws.Range("L116").Value = 1302.5714
$ws.Range("M116").Value = -66.63650000000007
$ws.Range("N116").Value = -5890.5714
$ws.Range("H122").Value = 4083.037
$ws.Range("I122").Value = 3932.3845
$ws.Range("K122").Value = 11797.1535
$ws.Range("M122").Value = -9347.1535
$ws.Range("H132").Value = 4529.1763
$ws.Range("I132").Value = 3399.7334
$ws.Range("J132").Value = 13000
$ws.Range("K132").Value = 10199.2002
$ws.Range("L132").Value = 39000
$ws.Range("M132").Value = -7669.200199999999
$ws.Range("N132").Value = -44060
$ws.Range("H136").Value = 4687.048
$ws.Range("I136").Value = 1732.1177
$ws.Range("K136").Value = 5196.3531
$ws.Range("M136").Value = -2646.3531
$ws.Range("N107").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1949.1666
$ws.Range("I3").Value = 2360.6365
$ws.Range("J3").Value = 1302.5714
$ws.Range("K3").Value = 2360.6365
$ws.Range("L3").Value = 1302.5714
$ws.Range("M3").Value = -2246.6365
$ws.Range("N3").Value = -1530.5714
$ws.Range("H86").Value = 2007.2
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("H89").Value = 2007.2
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("H94").Value = 1032
$ws.Range("I94").Value = 1032
$ws.Range("K94").Value = 1032
$ws.Range("M94").Value = -581
$ws.Range("H99").Value = 2506
$ws.Range("I99").Value = 2194.25
$ws.Range("K99").Value = 2194.25
$ws.Range("M99").Value = -696.25
$ws.Range("H105").Value = 1471.8966
$ws.Range("I105").Value = 1446.1818
$ws.Range("K105").Value = 1446.1818
$ws.Range("M105").Value = 300.8181999999999
$ws.Range("H107").Value = 2080.9375
$ws.Range("I107").Value = 1745.3636
$ws.Range("K107").Value = 1745.3636
$ws.Range("M107").Value = 174.6364000000001
$ws.Range("N86").ClearContents()
$ws.Range("N89").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1356.9286
$ws.Range("J16").Value = 1931.2
$ws.Range("L16").Value = 1931.2
$ws.Range("N16").Value = -2505.2
$ws.Range("H31").Value = 8121.875
$ws.Range("J31").Value = 9670.117
$ws.Range("L31").Value = 9670.117
$ws.Range("N31").Value = -10260.117
$ws.Range("H34").Value = 8121.875
$ws.Range("J34").Value = 9670.117
$ws.Range("L34").Value = 9670.117
$ws.Range("N34").Value = -10074.117
$ws.Range("H50").Value = 43356.855
$ws.Range("I50").Value = 10001.5
$ws.Range("J50").Value = 56699
$ws.Range("K50").Value = 10001.5
$ws.Range("L50").Value = 56699
$ws.Range("M50").Value = -9376.5
$ws.Range("N50").Value = -57949
$ws.Range("H99").Value = 2100.5103
$ws.Range("I99").Value = 1912.8529
$ws.Range("J99").Value = 2525.8667
$ws.Range("K99").Value = 1912.8529
$ws.Range("L99").Value = 2525.8667
$ws.Range("M99").Value = -414.8529000000001
$ws.Range("N99").Value = -5521.8667
$ws.Range("H113").Value = 1356.9286
$ws.Range("J113").Value = 1931.2
$ws.Range("L113").Value = 1931.2
$ws.Range("N113").Value = -6271.2
$ws.Range("H126").Value = 2100.5103
$ws.Range("I126").Value = 1912.8529
$ws.Range("J126").Value = 2525.8667
$ws.Range("K126").Value = 5738.5587
$ws.Range("L126").Value = 7577.6001
$ws.Range("M126").Value = -3268.5587
$ws.Range("N126").Value = -12517.6001
$ws.Range("H132").Value = 4189.143
$ws.Range("I132").Value = 2939.5881
$ws.Range("J132").Value = 9499.75
$ws.Range("K132").Value = 8818.764299999999
$ws.Range("L132").Value = 28499.25
$ws.Range("M132").Value = -6288.764299999999
$ws.Range("N132").Value = -33559.25
$ws.Range("H134").Value = 5567.273
$ws.Range("I134").Value = 4474.15
$ws.Range("K134").Value = 13422.45
$ws.Range("M134").Value = -10887.45
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2608.375
$ws.Range("I3").Value = 2608.375
$ws.Range("K3").Value = 7825.125
$ws.Range("M3").Value = -7713.125
$ws.Range("H37").Value = 216666.67
$ws.Range("J37").Value = 216666.67
$ws.Range("L37").Value = 650000.01
$ws.Range("N37").Value = -650224.01
$ws.Range("H46").Value = 844.7778
$ws.Range("I46").Value = 281
$ws.Range("K46").Value = 843
$ws.Range("M46").Value = -752
$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 3000
$ws.Range("M125").Value = 1920
$ws.Range("H132").Value = 4051.8
$ws.Range("I132").Value = 4103.6
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 36932.4
$ws.Range("L132").Value = 36000
$ws.Range("M132").Value = -34402.4
$ws.Range("N132").Value = -41060
$ws.Range("H133").Value = 10485.1875
$ws.Range("I133").Value = 3724.625
$ws.Range("K133").Value = 11173.875
$ws.Range("M133").Value = -6113.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3934.0833
$ws.Range("I102").Value = 2213.75
$ws.Range("K102").Value = 2213.75
$ws.Range("M102").Value = -591.75
$ws.Range("H126").Value = 4411.2354
$ws.Range("I126").Value = 2691.6155
$ws.Range("K126").Value = 8074.8465
$ws.Range("M126").Value = -5604.8465

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2109.4285
$ws.Range("I16").Value = 2118.2307
$ws.Range("K16").Value = 2118.2307
$ws.Range("M16").Value = -1948.2307
$ws.Range("H46").Value = 3410.8572
$ws.Range("J46").Value = 3923.5
$ws.Range("L46").Value = 3923.5
$ws.Range("N46").Value = -4299.5
$ws.Range("H55").Value = 882.6842
$ws.Range("I55").Value = 212.625
$ws.Range("K55").Value = 212.625
$ws.Range("M55").Value = -39.625
$ws.Range("H93").Value = 2551.3635
$ws.Range("I93").Value = 2582.2856
$ws.Range("K93").Value = 2582.2856
$ws.Range("M93").Value = -1334.2856
$ws.Range("H100").Value = 7692.1113
$ws.Range("I100").Value = 7653.625
$ws.Range("K100").Value = 7653.625
$ws.Range("M100").Value = -7112.625

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 23866.666
$ws.Range("J112").Value = 23866.666
$ws.Range("L112").Value = 23866.666
$ws.Range("N112").Value = -26820.666
$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -84800
$ws.Range("H136").Value = 7032.1304
$ws.Range("J136").Value = 10986.75
$ws.Range("L136").Value = 32960.25
$ws.Range("N136").Value = -38060.25
